{"js": "// Replace the hyphen-style date ranges with en-dashes in the\n// \"2022: Datumi kampanje za opazovanje Herkulovo ozvezdje: ...\" sentence.\n// This exact sentence appears 4 times in the document; every occurrence\n// gets the same three hyphen -> en-dash (U+2013) substitutions.\n\nconst oldText = \"2022: Datumi kampanje za opazovanje Herkulovo ozvezdje: 13.-22. junij, 12.-21. julij, 10.-19. avgust\";\nconst newText = \"2022: Datumi kampanje za opazovanje Herkulovo ozvezdje: 13.\\u201322. junij, 12.\\u201321. julij, 10.\\u201319. avgust\";\n\nconst body = context.document.body;\nconst results = body.search(oldText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newText, \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Replace the hyphen-style date ranges with en-dashes in the\n# \"2022: Datumi kampanje za opazovanje Herkulovo ozvezdje: ...\" sentence.\n# The same sentence (with plain hyphens) occurs 4 times in the document;\n# every occurrence gets the same three \"-\" -> \"\u2013\" (en dash) substitutions,\n# so a single document-wide Find/Replace-All covers all of them.\n\n$d = $word.ActiveDocument\n\n$findText    = \"13.-22. junij, 12.-21. julij, 10.-19. avgust\"\n$replaceText = \"13.\u201322. junij, 12.\u201321. julij, 10.\u201319. avgust\"\n\n$rng = $d.Content\n$rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n"}
